$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.272.73'
$ws.Range("E2").Value = '  +5.38%  '
$ws.Range("D3").Value = '1.916.78'
$ws.Range("E3").Value = '  +5.75%  '
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '253.86'
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.5145'
$ws.Range("E7").Value = '  +3.14%  '
$ws.Range("D8").Value = '46.01'
$ws.Range("E8").Value = '  +6.69%  '
$ws.Range("D9").Value = '0.2987'
$ws.Range("E9").Value = '  +7.47%  '
$ws.Range("D10").Value = '0.06831'
$ws.Range("E10").Value = '  +6.89%  '
$ws.Range("D11").Value = '1.915.47'
$ws.Range("E11").Value = '  +5.74%  '
$ws.Range("D12").Value = '17.46'
$ws.Range("E12").Value = '  +4.24%  '
$ws.Range("D13").Value = '0.07358'
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '0.6958'
$ws.Range("E14").Value = '  +7.39%  '
$ws.Range("D15").Value = '87.87'
$ws.Range("E15").Value = '  +7.36%  '
$ws.Range("D16").Value = '4.908'
$ws.Range("E16").Value = '  +4.39%  '
$ws.Range("D17").Value = '30.266.46'
$ws.Range("E17").Value = '  +5.46%  '
$ws.Range("D18").Value = '0.000007975'
$ws.Range("E18").Value = '  +7.88%  '
$ws.Range("D19").Value = '0.9992'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '13.07'
$ws.Range("E20").Value = '  +6.53%  '
$ws.Range("D21").Value = '2.162.29'
$ws.Range("E21").Value = '  +6.02%  '
$ws.Range("D22").Value = '0.9981'
$ws.Range("D23").Value = '4.863'
$ws.Range("E23").Value = '  +5.23%  '
$ws.Range("D24").Value = '5.744'
$ws.Range("E24").Value = '  +7.72%  '
$ws.Range("E25").Value = '  +3.37%  '
$ws.Range("D26").Value = '146.39'
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("D27").Value = '138.81'
$ws.Range("E27").Value = '  +23.14%  '
$ws.Range("D28").Value = '17.31'
$ws.Range("E28").Value = '  +8.14%  '
$ws.Range("D29").Value = '2.026'
$ws.Range("E29").Value = '  +8.03%  '
$ws.Range("D30").Value = '1.385'
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("D31").Value = '4.281'
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("D32").Value = '0.08848'
$ws.Range("E32").Value = '  +5.96%  '
$ws.Range("D33").Value = '4.032'
$ws.Range("E33").Value = '  +5.10%  '
$ws.Range("D34").Value = '0.05138'
$ws.Range("E34").Value = '  +3.19%  '
$ws.Range("D35").Value = '1.164'
$ws.Range("E35").Value = '  +6.87%  '
$ws.Range("D36").Value = '0.7186'
$ws.Range("E36").Value = '  +6.34%  '
$ws.Range("D37").Value = '2.685'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("D38").Value = '2.847'
$ws.Range("E38").Value = '  +4.40%  '
$ws.Range("D39").Value = '2.318'
$ws.Range("E39").Value = '  +7.89%  '
$ws.Range("D40").Value = '0.9777'
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("D41").Value = '0.01702'
$ws.Range("E41").Value = '  +6.99%  '
$ws.Range("D42").Value = '6.106'
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("D43").Value = '106.45'
$ws.Range("E43").Value = '  +4.95%  '
$ws.Range("D44").Value = '0.4329'
$ws.Range("E44").Value = '  +5.26%  '
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '7.739'
$ws.Range("E46").Value = '  +7.56%  '
$ws.Range("D47").Value = '0.1287'
$ws.Range("E47").Value = '  +5.13%  '
$ws.Range("D48").Value = '0.05755'
$ws.Range("E48").Value = '  +4.69%  '
$ws.Range("D49").Value = '33.54'
$ws.Range("E49").Value = '  +6.59%  '
$ws.Range("D50").Value = '8.558'
$ws.Range("E50").Value = '  +4.46%  '
$ws.Range("D51").Value = '0.3845'
$ws.Range("E51").Value = '  +5.88%  '
